# Trading update: 2026-02-17 20:23:43
# Appends 8 new MarketMaking trade rows (trade #36 .. #43) to the
# "All Trades" sheet (starting at row 37) and to the "MarketMaking"
# strategy sheet (starting at row 4), mirroring the original workbook's
# per-trade logging behaviour.

$wb = $excel.ActiveWorkbook

# Each element: Trade#, Date, Time, Strategy, Side, EntryPrice, Status,
#               PnLPct, PnLDollar, CapitalAfter, EntrySlippage, ExitSlippage,
#               Confidence, EntryReason, Duration
$newTrades = @(
    @(36, "2026-02-17", "20:22:31", "MarketMaking", "UP",   0.86, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(37, "2026-02-17", "20:22:38", "MarketMaking", "UP",   0.85, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(38, "2026-02-17", "20:22:45", "MarketMaking", "DOWN", 0.13, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(39, "2026-02-17", "20:22:52", "MarketMaking", "UP",   0.86, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(40, "2026-02-17", "20:23:11", "MarketMaking", "DOWN", 0.14, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(41, "2026-02-17", "20:23:18", "MarketMaking", "DOWN", 0.17, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(42, "2026-02-17", "20:23:25", "MarketMaking", "UP",   0.84, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(43, "2026-02-17", "20:23:38", "MarketMaking", "UP",   0.86, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0)
)

function Write-TradeRow($ws, $row, $trade) {
    $ws.Cells.Item($row, 1).Value = $trade[0]

    # Date / Time are stored as plain literal text in this log (not real
    # Excel dates), so force Text format while writing then drop back to
    # the default "Normal" style to avoid leaving a stray number format
    # behind on the cell.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $trade[1]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $trade[2]
    $ws.Cells.Item($row, 4).Value = $trade[3]
    $ws.Cells.Item($row, 5).Value = $trade[4]
    $ws.Cells.Item($row, 6).Value = $trade[5]
    # Column G (Exit Price) left blank - trade is still OPEN.
    $ws.Cells.Item($row, 8).Value = $trade[6]
    $ws.Cells.Item($row, 9).Value = $trade[7]
    $ws.Cells.Item($row, 10).Value = $trade[8]
    $ws.Cells.Item($row, 11).Value = $trade[9]
    $ws.Cells.Item($row, 12).Value = $trade[10]
    $ws.Cells.Item($row, 13).Value = $trade[11]
    $ws.Cells.Item($row, 14).Value = $trade[12]
    $ws.Cells.Item($row, 15).Value = $trade[13]
    # Column P (Exit Reason) left blank - trade is still OPEN.
    $ws.Cells.Item($row, 17).Value = $trade[14]
}

# "All Trades" sheet already has 36 data rows (header + 35 trades, trade
# #35 on row 36) -> new trades continue on row 37.
$wsAll = $wb.Worksheets.Item("All Trades")
$startRowAll = 37
for ($i = 0; $i -lt $newTrades.Count; $i++) {
    $targetRow = $startRowAll + $i
    $trade = $newTrades[$i]
    Write-TradeRow $wsAll $targetRow $trade
}

# "MarketMaking" sheet only tracks this strategy's trades and already has
# 2 data rows (trade #34 on row 2, #35 on row 3) -> continues on row 4.
$wsStrategy = $wb.Worksheets.Item("MarketMaking")
$startRowStrategy = 4
for ($i = 0; $i -lt $newTrades.Count; $i++) {
    $targetRow = $startRowStrategy + $i
    $trade = $newTrades[$i]
    Write-TradeRow $wsStrategy $targetRow $trade
}
